# Continent_Country.xlsx - add two new Oceania entries
# ("Wallis and Futuna" under France, and "Micronesia") as two new rows
# inserted right after the last existing Oceania row (old row 254, the
# "Samoa" entry) and before the South America block (old row 255). This
# mirrors how the sheet had previously grown: new territories appended to
# the bottom of the Oceania block rather than a full re-sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 255-256, pushing the South America block (and
# everything after it) down by two rows. Excel copies the formatting of
# the row immediately above the insertion point (the "Samoa" row, which
# uses the wrap-text / vertical-center style used for territory rows),
# so the new B/C cells automatically pick up the right look.
$ws.Rows("255:256").Insert()

$ws.Range("A255").Value2 = "Oceania"
$ws.Range("B255").Value2 = "France"
$ws.Range("C255").Value2 = "Wallis and Futuna"

$ws.Range("A256").Value2 = "Oceania"
$ws.Range("B256").Value2 = "Micronesia"
$ws.Range("C256").Value2 = "Micronesia"

# The old row 260 (Costa Rica, A260 carried a highlighted style) is now
# row 262; drop that leftover cell highlight there so only the row that
# originally carried it two rows up (old row 262, the France / French
# Guiana entry, now row 264) keeps it.
$ws.Range("A262").ClearFormats()

# Reflect the user's final selection / scroll position while reviewing
# the newly-added rows.
$ws.Range("A256").Select()
$excel.ActiveWindow.ScrollRow = 226
$excel.ActiveWindow.ScrollColumn = 1
